# Auto update Excel log
# Appends new sensor-log rows to the ALERTS, PIR, Humidity and mmWave
# sheets, matching the source system's latest export.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($ws, $row, $col, $val)
    # Route every write through an explicit Text number format so values
    # that look like dates/times/percentages (e.g. "2026-01-30", "86.8%")
    # are stored as literal text instead of being auto-converted by Excel's
    # value parser. ClearFormats() afterwards drops the temporary format so
    # the cell is left with the workbook's default (unstyled) appearance.
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

function Add-LogRows {
    param($ws, $rows)
    foreach ($r in $rows) {
        $rowNum = $r[0]
        Set-TextCell $ws $rowNum 1 $r[1]
        Set-TextCell $ws $rowNum 2 $r[2]
        Set-TextCell $ws $rowNum 3 $r[3]
        Set-TextCell $ws $rowNum 4 $r[4]
        Set-TextCell $ws $rowNum 5 $r[5]
        Set-TextCell $ws $rowNum 6 $r[6]
    }
}

# ---- ALERTS: one new CRITICAL fall-detection alert ----
$wsAlerts = $wb.Worksheets.Item("ALERTS")
$alertsRows = @()
$alertsRows += ,@(11, "2026-01-30", "15:51:55", "15:00", "Living Room", "CRITICAL", "FALL_DETECTED")
Add-LogRows $wsAlerts $alertsRows

# ---- PIR: fourteen new motion-sensor readings ----
$wsPir = $wb.Worksheets.Item("PIR")
$pirRows = @()
$pirRows += ,@(167, "2026-01-30", "15:51:38", "15:00", "Bathroom", "No Motion", "Inactive")
$pirRows += ,@(168, "2026-01-30", "15:51:41", "15:00", "Bathroom", "No Motion", "Inactive")
$pirRows += ,@(169, "2026-01-30", "15:51:46", "15:00", "Bathroom", "No Motion", "Inactive")
$pirRows += ,@(170, "2026-01-30", "15:51:51", "15:00", "Bathroom", "No Motion", "Inactive")
$pirRows += ,@(171, "2026-01-30", "15:51:56", "15:00", "Bathroom", "No Motion", "Inactive")
$pirRows += ,@(172, "2026-01-30", "15:51:58", "15:00", "Living Room", "RECOVERY_DETECTION", "Inactive")
$pirRows += ,@(173, "2026-01-30", "15:52:01", "15:00", "Bathroom", "No Motion", "Inactive")
$pirRows += ,@(174, "2026-01-30", "15:52:06", "15:00", "Bathroom", "No Motion", "Inactive")
$pirRows += ,@(175, "2026-01-30", "15:52:11", "15:00", "Bathroom", "No Motion", "Inactive")
$pirRows += ,@(176, "2026-01-30", "15:52:16", "15:00", "Bathroom", "No Motion", "Inactive")
$pirRows += ,@(177, "2026-01-30", "15:52:21", "15:00", "Bathroom", "No Motion", "Inactive")
$pirRows += ,@(178, "2026-01-30", "15:52:26", "15:00", "Bathroom", "No Motion", "Inactive")
$pirRows += ,@(179, "2026-01-30", "15:52:31", "15:00", "Bathroom", "No Motion", "Inactive")
$pirRows += ,@(180, "2026-01-30", "15:52:36", "15:00", "Bathroom", "No Motion", "Inactive")
Add-LogRows $wsPir $pirRows

# ---- Humidity: five new readings ----
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @()
$humidityRows += ,@(98, "2026-01-30", "15:51:46", "15:00", "Bathroom", "86.8%", "Active")
$humidityRows += ,@(99, "2026-01-30", "15:51:55", "15:00", "Bathroom", "87.7%", "Active")
$humidityRows += ,@(100, "2026-01-30", "15:51:56", "15:00", "Bathroom", "86.8%", "Active")
$humidityRows += ,@(101, "2026-01-30", "15:52:11", "15:00", "Bathroom", "87.8%", "Active")
$humidityRows += ,@(102, "2026-01-30", "15:52:31", "15:00", "Bathroom", "87.8%", "Active")
Add-LogRows $wsHumidity $humidityRows

# ---- mmWave: five new presence-detection readings ----
$wsMmwave = $wb.Worksheets.Item("mmWave")
$mmwaveRows = @()
$mmwaveRows += ,@(19, "2026-01-30", "15:51:38", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")
$mmwaveRows += ,@(20, "2026-01-30", "15:51:58", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")
$mmwaveRows += ,@(21, "2026-01-30", "15:52:08", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")
$mmwaveRows += ,@(22, "2026-01-30", "15:52:19", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")
$mmwaveRows += ,@(23, "2026-01-30", "15:52:29", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")
Add-LogRows $wsMmwave $mmwaveRows
